# Update the crypto price list with the latest scraped values.
# GitHub Actions scheduled refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.503.59"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.427.99"
$ws.Range("E3").Value = "  +0.93%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.34"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.65"
$ws.Range("E6").Value = "  -2.40%  "

# Row 7 - now LidoStakedEther (was USDC) -- rows 7/8 swapped order
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.421.97"
$ws.Range("E7").Value = "  +0.81%  "

# Row 8 - now USDC (was LidoStakedEther)
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  +1.39%  "

# Row 11 - Cardano
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.583"
$ws.Range("E11").Value = "  -1.25%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.69"
$ws.Range("E12").Value = "  -0.47%  "

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000281"
$ws.Range("E13").Value = "  -1.39%  "

# Row 14 - BitcoinCash
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "692.17"
$ws.Range("E14").Value = "  +1.05%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.975.93"
$ws.Range("E15").Value = "  +0.84%  "

# Row 16 - Polkadot
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.62"
$ws.Range("E16").Value = "  -0.40%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.551.50"
$ws.Range("E17").Value = "  +0.18%  "

# Row 18 - now TRON (was WrappedEther) -- rows 18/19 swapped order
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.121"
$ws.Range("E18").Value = "  +1.01%  "

# Row 19 - now WrappedEther (was TRON)
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.418.23"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.65"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("E21").Value = "  -0.51%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.895"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23 - Toncoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.41"
$ws.Range("E23").Value = "  -0.57%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.93"
$ws.Range("E24").Value = "  -1.12%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.28"
$ws.Range("E25").Value = "  -2.44%  "

# Row 26 - PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").Value = "  -0.96%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  -2.99%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.56"
$ws.Range("E28").Value = "  -0.78%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.52"
$ws.Range("E29").Value = "  -2.28%  "

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.75"
$ws.Range("E30").Value = "  +0.14%  "

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.07"
$ws.Range("E31").Value = "  +0.91%  "

# Row 32 - Bittensor
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "573.88"
$ws.Range("E32").Value = "  +3.17%  "

# Row 33 - dogwifhat
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.71"
$ws.Range("E33").Value = "  +0.62%  "

# Row 34 - Cosmos
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.01"
$ws.Range("E34").Value = "  -1.85%  "

# Row 35 - OKB
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.26"
$ws.Range("E35").Value = "  +0.23%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  -2.87%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  -0.04%  "

# Row 38 - Maker
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.583.34"
$ws.Range("E38").Value = "  -3.27%  "

# Row 39 - Kaspa
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  -1.08%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.88"
$ws.Range("E40").Value = "  -0.70%  "

# Row 41 - PEPE
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0732"
$ws.Range("E41").Value = "  +3.74%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 - Fetch.AI
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44 - ApeXProtocol
$ws.Range("E44").Value = "  +3.97%  "

# Row 45 - TheGraph
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").Value = "  -2.38%  "

# Row 46 - VeChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0417"
$ws.Range("E46").Value = "  -1.55%  "

# Row 47 - Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.45"
$ws.Range("E47").Value = "  +4.60%  "

# Row 48 - ThetaToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  -0.37%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -1.27%  "

# Row 50 - FirstDigitalUSD
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.17%  "

# Row 51 - Monero
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.62"
$ws.Range("E51").Value = "  +0.09%  "
